$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Recorded By" text cleanups (strip stray "system,"/"System," prefix) ---
$ws.Range("G2").Value  = "backup@backdoor.com"
$ws.Range("G6").Value  = "dnasr281@gmail.com"
$ws.Range("G12").Value = "dnasr281@gmail.com"
$ws.Range("G13").Value = "dnasr281@gmail.com"
$ws.Range("G29").Value = "backup@backdoor.com"
$ws.Range("G33").Value = "dnasr281@gmail.com"
$ws.Range("G39").Value = "dnasr281@gmail.com"
$ws.Range("G40").Value = "dnasr281@gmail.com"
$ws.Range("G56").Value = "backup@backdoor.com"
$ws.Range("G60").Value = "dnasr281@gmail.com"
$ws.Range("G66").Value = "dnasr281@gmail.com"
$ws.Range("G67").Value = "dnasr281@gmail.com"

# --- Reorder "Recorded By" lists ---
$ws.Range("G90").Value  = "dnasr281@gmail.com, admin@admin.com"
$ws.Range("G116").Value = "dnasr281@gmail.com, admin@admin.com"
$ws.Range("G142").Value = "dnasr281@gmail.com, admin@admin.com"

# --- Attendance count corrections ---
$ws.Range("H2").Value   = "32/53"
$ws.Range("H6").Value   = "43/53"
$ws.Range("H12").Value  = "30/53"
$ws.Range("H13").Value  = "35/53"
$ws.Range("H92").Value  = "45/56"
$ws.Range("H118").Value = "46/55"

# --- Class statistics ---
# Percentages are stored as literal text (not numeric %) in this sheet. A
# leading apostrophe stops Excel's automatic "looks like a percent" number
# conversion, but it also tags the cell's style with a quote-prefix flag, so
# afterwards a formats-only paste from an unrelated, still-plain cell that
# already carries the same base style (style 4 here, via L8) restores the
# original style index while leaving the freshly-written text value intact.
function Set-PlainText($a1, $text) {
    $ws.Range($a1).Value = "'" + $text
    $ws.Range("L8").Copy()
    $ws.Range($a1).PasteSpecial(-4122)
}

$ws.Range("L6").Value  = 65
$ws.Range("L7").Value  = 4
Set-PlainText "L9"  "40.9%"
Set-PlainText "L10" "61.4%"
Set-PlainText "S15" "59.4%"

# --- Group statistics rows 18-20 ---
$ws.Range("O18").Value = 10
$ws.Range("P18").Value = 1
Set-PlainText "R18" "38.5%"
Set-PlainText "S18" "63.0%"

$ws.Range("O19").Value = 10
$ws.Range("P19").Value = 1
Set-PlainText "R19" "38.5%"
Set-PlainText "S19" "67.3%"

$ws.Range("O20").Value = 10
$ws.Range("P20").Value = 1
Set-PlainText "R20" "38.5%"
Set-PlainText "S20" "70.4%"

# --- Rows 93 / 119 / 145: sessions flip from Recorded -> Not Recorded ---
# These rows need both a value change (the session becomes unrecorded) and a
# format change (green "Recorded" shading -> pink "Not Recorded" shading,
# matching the existing style used by row 31's "Not Recorded" entry).
#
# Order matters: write the new cell values FIRST (the D/E text values and the
# blank "Recorded By" cell use a leading apostrophe so Excel keeps them as
# literal text instead of re-parsing "11"/"01/10/2025" as a number/date or
# collapsing an empty string to a blank cell), THEN copy row 31's formatting
# onto the row with a formats-only paste so the values just written are kept.
function Set-NotRecorded($rowNum, $group, $total) {
    $ws.Cells.Item($rowNum, 1).Value = "Year 4"
    $ws.Cells.Item($rowNum, 2).Value = $group
    $ws.Cells.Item($rowNum, 3).Value = "GENERAL SURGERY"
    $ws.Cells.Item($rowNum, 4).Value = "'11"
    $ws.Cells.Item($rowNum, 5).Value = "'01/10/2025"
    $ws.Cells.Item($rowNum, 6).Value = "10:30:00"
    $ws.Cells.Item($rowNum, 7).Value = "'"
    $ws.Cells.Item($rowNum, 8).Value = "0/" + $total
    $ws.Cells.Item($rowNum, 9).Value = "Not Recorded"

    $src = $ws.Range("A31:I31")
    $src.Copy()
    $dst = $ws.Range("A" + $rowNum + ":I" + $rowNum)
    $dst.PasteSpecial(-4122)
}

Set-NotRecorded 93  "B2D" "56"
Set-NotRecorded 119 "B2E" "55"
Set-NotRecorded 145 "B2F" "57"
